$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were added for "Camote" at the top of the
# historical data block (current rows 78-168 all shift down by two rows,
# ending at rows 80-170). Insert two blank rows at row 78 to produce that
# shift while preserving the existing row formatting (e.g. the date style
# on column D).
$ws.Rows("78:79").Insert()

# Populate the first new record (row 78).
$ws.Range("A78").Value = 9
$ws.Range("B78").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C78").Value = "Metropolitana"
$ws.Range("D78").Value = 45117
$ws.Range("E78").Value = 13
$ws.Range("F78").Value = 100114002
$ws.Range("G78").Value = "Camote"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 430
$ws.Range("K78").Value = 19000
$ws.Range("L78").Value = 20000
$ws.Range("M78").Value = 19500
$ws.Range("N78").Value = "$/caja 18 kilos"
$ws.Range("O78").Value = "Perú"
$ws.Range("P78").Value = 1083
$ws.Range("Q78").Value = 18
$ws.Range("R78").Value = "Hortaliza"

# Populate the second new record (row 79).
$ws.Range("A79").Value = 9
$ws.Range("B79").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C79").Value = "Metropolitana"
$ws.Range("D79").Value = 45117
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = 100114002
$ws.Range("G79").Value = "Camote"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 700
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 16000
$ws.Range("M79").Value = 15500
$ws.Range("N79").Value = "$/malla 18 kilos"
$ws.Range("O79").Value = "Perú"
$ws.Range("P79").Value = 861
$ws.Range("Q79").Value = 18
$ws.Range("R79").Value = "Hortaliza"
